$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matrix")

# Core data change: raise the "red" threshold (J10) from 7 to 8.
# This single change cascades through every dependent formula in the sheet
# (column C IF() results, column D sums, and the H12/H13 helper totals).
$ws.Range("J10").Value = 8

# New label cell introduced alongside the threshold change.
$ws.Range("H8").Value = "limit"

# Selection/view state update recorded in the sheet view.
$ws.Range("A1:D1048576").Select()
